$d = $word.ActiveDocument

$replacements = @(
    @{old="773×7="; new="883×5="},
    @{old="776×2="; new="195×8="},
    @{old="916×4="; new="883×3="},
    @{old="738×9="; new="435×9="},
    @{old="762×4="; new="525×2="},
    @{old="964×2="; new="961×5="},
    @{old="432×8="; new="767×3="},
    @{old="184×9="; new="872×8="},
    @{old="504×8="; new="729×8="},
    @{old="624×5="; new="251×4="},
    @{old="695×4="; new="519×4="},
    @{old="370×7="; new="547×7="},
    @{old="314×9="; new="811×4="},
    @{old="780×8="; new="516×6="},
    @{old="406×4="; new="517×4="},
    @{old="367×7="; new="519×6="},
    @{old="278×8="; new="563×4="},
    @{old="147×6="; new="621×9="},
    @{old="619×9="; new="101×8="},
    @{old="563×3="; new="224×9="},
    @{old="362×9="; new="722×4="},
    @{old="814×3="; new="228×3="},
    @{old="603×2="; new="719×5="},
    @{old="396×2="; new="122×3="},
    @{old="631×3="; new="715×3="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
